# Insert a new row above row 50 (shifts existing rows 50..173 down to 51..174)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record's data
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 44987
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 100112052
$ws.Range("G50").Value = "Albahaca"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 70
$ws.Range("K50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = 6000
$ws.Range("N50").Value = "$/docena de matas"
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 1000
$ws.Range("Q50").Value = 6
$ws.Range("R50").Value = "Hortaliza"
